# Add small slide-number text boxes ("1".."5") to the bottom-left corner
# of slides 2-6 (the content slides; the title slide and the closing
# "Thank You" slide are left untouched).
#
# Each new text box:
#   - is an auto-sized, non-wrapping TextBox (wrap="none" + spAutoFit)
#   - has no fill
#   - sits at off (x=0, y=6488668 EMU) with size (cx=300082, cy=369332) EMU
#   - contains a single run with the slide's number as text
#
# PowerPoint's COM object model expresses shape geometry in points, while
# the OOXML stores EMU (1 pt = 12700 EMU), so we convert before calling
# AddTextbox.

$p = $ppt.ActivePresentation

$EMU_PER_POINT = 12700

$offX  = 0 / $EMU_PER_POINT
$offY  = 6488668 / $EMU_PER_POINT
$width = 300082 / $EMU_PER_POINT
$height= 369332 / $EMU_PER_POINT

$msoFalse = 0
$msoTextOrientationHorizontal = 1
$ppAutoSizeShapeToFitText = 1

$slideNumbers = @{
    2 = "1"
    3 = "2"
    4 = "3"
    5 = "4"
    6 = "5"
}

foreach ($slideIndex in 2..6) {
    $s = $p.Slides.Item($slideIndex)
    $shp = $s.Shapes.AddTextbox($msoTextOrientationHorizontal, $offX, $offY, $width, $height)
    $shp.TextFrame.WordWrap = $msoFalse
    $shp.TextFrame.AutoSize = $ppAutoSizeShapeToFitText
    $shp.Fill.Visible = $msoFalse
    $shp.TextFrame.TextRange.Text = $slideNumbers[$slideIndex]
}
